# edit.ps1 - PowerPoint COM-interop script
#
# Reproduces the commit:
#   1. Three tables (on the slides that hold the "Business structure" /
#      comparison tables) get their table style switched from the custom
#      "Table_0" style ({30BDE7A2-EF52-4911-B817-F1FB23B921EE}) to the
#      built-in "Medium Style 2 - Accent 1" style
#      ({CB26E7FF-3F0C-4B48-8923-98698A645FAA}).
#   2. The deck's theme palette is switched from the "Integral" / "Red
#      Violet" colour scheme back to the default Office colour scheme
#      (the presentation-level theme part that everything except the
#      notes master is rendered with).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Table styles -----------------------------------------------------------
# ---------------------------------------------------------------------------
# Every slide that contains a table in this deck has exactly one table
# shape; walk the whole deck instead of hard-coding slide numbers so the
# script is robust to the exact slide indices.

$newTableStyle = "{CB26E7FF-3F0C-4B48-8923-98698A645FAA}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Theme colours -----------------------------------------------------------
# ---------------------------------------------------------------------------
# Restore the 12 standard "Office" theme colours (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) on the presentation's active theme.

function HexToRgbVal([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = HexToRgbVal($officeColors[$i - 1])
}
